# Outer-join cleanup: drop the rows whose joined columns came back NaN
# (missing studio / missing imdb_rating). In the original sheet these are
# movie_id 110 (Bajirao Mastani, row 10), 124 (Parasite, row 22),
# 128 (Taare Zameen Par, row 26) and 131 (Sanju, row 29).
# Delete from the bottom up so earlier row numbers stay valid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(29).Delete()
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(10).Delete()
